$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"

$dCell = $ws.Cells.Item($row, 4)
$dCell.Value = 44446
$dCell.NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = "Otros"
$ws.Cells.Item($row, 9).Value = 100107002
$ws.Cells.Item($row, 10).Value = "Chirimoya"
$ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 3200
$ws.Cells.Item($row, 15).Value = 3300
$ws.Cells.Item($row, 16).Value = 3250
$ws.Cells.Item($row, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 19).Value = 3250
$ws.Cells.Item($row, 20).Value = 1
